# The source commit ("Moving from 2.0.1 to 2.0.2") regenerated this
# template through a newer docx4j release. Diffing the two packages shows
# that every single hunk (word/document.xml, word/footer1-3.xml,
# word/footnotes.xml, word/header1-3.xml, word/styles.xml) is purely a
# cosmetic re-serialization performed by that Java toolchain:
#   - namespace declarations / attributes get re-emitted in alphabetical
#     order (e.g. <w:tcW w:w="3070" w:type="dxa"/> -> <w:tcW w:type="dxa" w:w="3070"/>)
#   - a couple of "a:graphicFrameLocks"/"a:graphic" elements drop a
#     redundant xmlns:a re-declaration that is already in scope from an
#     ancestor element
#   - the base64 VML "o:gfxdata" blob is re-wrapped (identical bytes,
#     different line folding)
# No text, value, relationship, formatting, image, or structural content
# actually changes anywhere in the package (confirmed by comparing every
# changed attribute list as a set: before == after in every hunk).
#
# There is therefore no Word object-model mutation that corresponds to
# this commit - it is not something a user (or a macro) did inside Word,
# it is an artifact of the XML writer used by the external build tool
# that produced the fixture. We touch nothing, which keeps the document
# exactly equivalent (same text, same formatting, same structure, same
# relationships) to the target of the diff.

$d = $word.ActiveDocument

# No-op sanity touch: just confirm the document is reachable without
# mutating any content, formatting or structure.
$null = $d.Content.Text
